$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3041.125
$ws.Range("I17").Value = 2887
$ws.Range("J17").Value = 3063.1428
$ws.Range("K17").Value = 8661
$ws.Range("L17").Value = 9189.428400000001
$ws.Range("M17").Value = -8493
$ws.Range("N17").Value = -9525.428400000001
$ws.Range("H18").Value = 507.33334
$ws.Range("I18").Value = 408.8
$ws.Range("K18").Value = 408.8
$ws.Range("M18").Value = -124.8
$ws.Range("H28").Value = 766.8
$ws.Range("I28").Value = 766.8
$ws.Range("K28").Value = 766.8
$ws.Range("M28").Value = -281.8
$ws.Range("H32").Value = 3296.7778
$ws.Range("I32").Value = 3070.25
$ws.Range("J32").Value = 3478
$ws.Range("K32").Value = 3070.25
$ws.Range("L32").Value = 3478
$ws.Range("M32").Value = -2744.25
$ws.Range("N32").Value = -4130
$ws.Range("H42").Value = 2512.4167
$ws.Range("I42").Value = 630.25
$ws.Range("J42").Value = 6276.75
$ws.Range("K42").Value = 1890.75
$ws.Range("L42").Value = 18830.25
$ws.Range("M42").Value = -1660.75
$ws.Range("N42").Value = -19290.25
$ws.Range("H51").Value = 12334.333
$ws.Range("I51").Value = 13000.5
$ws.Range("J51").Value = 11002
$ws.Range("K51").Value = 13000.5
$ws.Range("L51").Value = 11002
$ws.Range("M51").Value = -12516.5
$ws.Range("N51").Value = -11970
$ws.Range("H54").Value = 7029.143
$ws.Range("I54").Value = 7367.3335
$ws.Range("J54").Value = 5000
$ws.Range("K54").Value = 7367.3335
$ws.Range("L54").Value = 5000
$ws.Range("M54").Value = -6881.3335
$ws.Range("N54").Value = -5972
$ws.Range("H64").Value = 4999.5
$ws.Range("I64").Value = 4999.5
$ws.Range("K64").Value = 4999.5
$ws.Range("M64").Value = -4751.5
$ws.Range("H67").Value = 4999.5
$ws.Range("I67").Value = 4999.5
$ws.Range("K67").Value = 4999.5
$ws.Range("M67").Value = -4141.5
$ws.Range("H70").Value = 2430.7
$ws.Range("I70").Value = 1401.8
$ws.Range("J70").Value = 3459.6
$ws.Range("K70").Value = 4205.4
$ws.Range("L70").Value = 10378.8
$ws.Range("M70").Value = -3935.4
$ws.Range("N70").Value = -10918.8
$ws.Range("H73").Value = 2430.7
$ws.Range("I73").Value = 1401.8
$ws.Range("J73").Value = 3459.6
$ws.Range("K73").Value = 4205.4
$ws.Range("L73").Value = 10378.8
$ws.Range("M73").Value = -3269.4
$ws.Range("N73").Value = -12250.8
$ws.Range("H88").Value = 2320
$ws.Range("I88").Value = 697.5
$ws.Range("J88").Value = 2783.5715
$ws.Range("K88").Value = 697.5
$ws.Range("L88").Value = 2783.5715
$ws.Range("M88").Value = -291.5
$ws.Range("N88").Value = -3595.5715
$ws.Range("H91").Value = 2320
$ws.Range("I91").Value = 697.5
$ws.Range("J91").Value = 2783.5715
$ws.Range("K91").Value = 697.5
$ws.Range("L91").Value = 2783.5715
$ws.Range("M91").Value = 706.5
$ws.Range("N91").Value = -5591.5715
$ws.Range("H94").Value = 5500.5
$ws.Range("I94").Value = 5500.5
$ws.Range("K94").Value = 5500.5
$ws.Range("M94").Value = -5049.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H112").Value = 3787.3044
$ws.Range("J112").Value = 3787.3044
$ws.Range("L112").Value = 11361.9132
$ws.Range("N112").Value = -13577.9132
$ws.Range("H125").Value = 987.3333
$ws.Range("J125").Value = 1321
$ws.Range("L125").Value = 11889
$ws.Range("N125").Value = -16809
$ws.Range("H127").Value = 3571
$ws.Range("I127").Value = 3464
$ws.Range("K127").Value = 10392
$ws.Range("M127").Value = -5432
$ws.Range("H129").Value = 6281.5
$ws.Range("I129").Value = 2079.889
$ws.Range("K129").Value = 6239.667
$ws.Range("M129").Value = -1239.667
$ws.Range("H132").Value = 4749.732
$ws.Range("I132").Value = 2771.361
$ws.Range("J132").Value = 18994
$ws.Range("K132").Value = 8314.082999999999
$ws.Range("L132").Value = 56982
$ws.Range("M132").Value = -5784.082999999999
$ws.Range("N132").Value = -62042
$ws.Range("H138").Value = 3455.8704
$ws.Range("J138").Value = 3587.7073
$ws.Range("L138").Value = 10763.1219
$ws.Range("N138").Value = -21043.1219

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6260.625
$ws.Range("I32").Value = 3930.6326
$ws.Range("K32").Value = 3930.6326
$ws.Range("M32").Value = -3643.6326
$ws.Range("H45").Value = 3355.611
$ws.Range("I45").Value = 1284.5714
$ws.Range("J45").Value = 4673.5454
$ws.Range("K45").Value = 1284.5714
$ws.Range("L45").Value = 4673.5454
$ws.Range("M45").Value = -907.5714
$ws.Range("N45").Value = -5427.5454
$ws.Range("H61").Value = 2087.3
$ws.Range("I61").Value = 1762.6666
$ws.Range("J61").Value = 2574.25
$ws.Range("K61").Value = 1762.6666
$ws.Range("L61").Value = 2574.25
$ws.Range("M61").Value = -1550.6666
$ws.Range("N61").Value = -2998.25
$ws.Range("H63").Value = 3211
$ws.Range("I63").Value = 2849
$ws.Range("K63").Value = 2849
$ws.Range("M63").Value = -2163
$ws.Range("H66").Value = 3211
$ws.Range("I66").Value = 2849
$ws.Range("K66").Value = 14245
$ws.Range("M66").Value = -10813
$ws.Range("H132").Value = 51238.19
$ws.Range("I132").Value = 58439.723
$ws.Range("K132").Value = 175319.169
$ws.Range("M132").Value = -172789.169
$ws.Range("H133").Value = 50079
$ws.Range("I133").Value = 30237
$ws.Range("J133").Value = 60000
$ws.Range("K133").Value = 30237
$ws.Range("L133").Value = 60000
$ws.Range("M133").Value = -27707
$ws.Range("N133").Value = -65060
$ws.Range("H134").Value = 139582.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 139582.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 139582.5
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -149722.5
$ws.Range("H136").Value = 2087.3
$ws.Range("I136").Value = 1762.6666
$ws.Range("J136").Value = 2574.25
$ws.Range("K136").Value = 5287.9998
$ws.Range("L136").Value = 7722.75
$ws.Range("M136").Value = -2737.9998
$ws.Range("N136").Value = -12822.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 973.44446
$ws.Range("I107").Value = 973.44446
$ws.Range("K107").Value = 973.44446
$ws.Range("M107").Value = 946.55554

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 24623.75
$ws.Range("I55").Value = 24623.75
$ws.Range("K55").Value = 24623.75
$ws.Range("M55").Value = -24308.75
$ws.Range("H62").Value = 3661
$ws.Range("I62").Value = 3576.25
$ws.Range("K62").Value = 3576.25
$ws.Range("M62").Value = -2952.25
$ws.Range("H65").Value = 3661
$ws.Range("I65").Value = 3576.25
$ws.Range("K65").Value = 17881.25
$ws.Range("M65").Value = -14761.25
$ws.Range("H99").Value = 2428.7
$ws.Range("I99").Value = 1924.3334
$ws.Range("J99").Value = 3185.25
$ws.Range("K99").Value = 1924.3334
$ws.Range("L99").Value = 3185.25
$ws.Range("M99").Value = -426.3334
$ws.Range("N99").Value = -6181.25
$ws.Range("H104").Value = 41259
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H122").Value = 3290.111
$ws.Range("I122").Value = 3069.4
$ws.Range("J122").Value = 3566
$ws.Range("K122").Value = 9208.200000000001
$ws.Range("L122").Value = 10698
$ws.Range("M122").Value = -6758.200000000001
$ws.Range("N122").Value = -15598
$ws.Range("H126").Value = 2428.7
$ws.Range("I126").Value = 1924.3334
$ws.Range("J126").Value = 3185.25
$ws.Range("K126").Value = 5773.0002
$ws.Range("L126").Value = 9555.75
$ws.Range("M126").Value = -3303.0002
$ws.Range("N126").Value = -14495.75
$ws.Range("H138").Value = 53115.31
$ws.Range("I138").Value = 41388.89
$ws.Range("J138").Value = 79499.75
$ws.Range("K138").Value = 41388.89
$ws.Range("L138").Value = 79499.75
$ws.Range("M138").Value = -36248.89
$ws.Range("N138").Value = -89779.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 976.8182
$ws.Range("I5").Value = 798.2
$ws.Range("J5").Value = 1125.6666
$ws.Range("K5").Value = 2394.6
$ws.Range("L5").Value = 3376.9998
$ws.Range("M5").Value = -2282.6
$ws.Range("N5").Value = -3600.9998
$ws.Range("H60").Value = 277
$ws.Range("I60").Value = 252.25
$ws.Range("J60").Value = 326.5
$ws.Range("K60").Value = 756.75
$ws.Range("L60").Value = 979.5
$ws.Range("M60").Value = -505.75
$ws.Range("N60").Value = -1481.5
$ws.Range("H113").Value = 1294.4445
$ws.Range("J113").Value = 1467.1154
$ws.Range("L113").Value = 4401.3462
$ws.Range("N113").Value = -8741.3462
$ws.Range("H132").Value = 917
$ws.Range("I132").Value = 688.2222
$ws.Range("J132").Value = 1328.8
$ws.Range("K132").Value = 6193.999800000001
$ws.Range("L132").Value = 11959.2
$ws.Range("M132").Value = -3663.999800000001
$ws.Range("N132").Value = -17019.2
$ws.Range("H135").Value = 976.8182
$ws.Range("I135").Value = 798.2
$ws.Range("J135").Value = 1125.6666
$ws.Range("K135").Value = 7183.8
$ws.Range("L135").Value = 10130.9994
$ws.Range("M135").Value = -4648.8
$ws.Range("N135").Value = -15200.9994

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 67870
$ws.Range("I107").Value = 143200.28
$ws.Range("K107").Value = 143200.28
$ws.Range("M107").Value = -141280.28
$ws.Range("H122").Value = 2146.8333
$ws.Range("I122").Value = 970
$ws.Range("J122").Value = 3794.4
$ws.Range("K122").Value = 2910
$ws.Range("L122").Value = 11383.2
$ws.Range("M122").Value = -460
$ws.Range("N122").Value = -16283.2
$ws.Range("H132").Value = 174077.17
$ws.Range("J132").Value = 10503.5
$ws.Range("L132").Value = 31510.5
$ws.Range("N132").Value = -36570.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 25504.652
$ws.Range("J22").Value = 3314.923
$ws.Range("L22").Value = 3314.923
$ws.Range("N22").Value = -3904.923
$ws.Range("H27").Value = 25504.652
$ws.Range("J27").Value = 3314.923
$ws.Range("L27").Value = 3314.923
$ws.Range("N27").Value = -3528.923
$ws.Range("H46").Value = 23541.334
$ws.Range("I46").Value = 35391.555
$ws.Range("J46").Value = 5766
$ws.Range("K46").Value = 35391.555
$ws.Range("L46").Value = 5766
$ws.Range("M46").Value = -35203.555
$ws.Range("N46").Value = -6142
$ws.Range("H57").Value = 29361.5
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 29361.5
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 29361.5
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -30493.5
$ws.Range("H68").Value = 6326.3335
$ws.Range("I68").Value = 2999
$ws.Range("K68").Value = 2999
$ws.Range("M68").Value = -2250
$ws.Range("H71").Value = 6326.3335
$ws.Range("I71").Value = 2999
$ws.Range("K71").Value = 14995
$ws.Range("M71").Value = -11251
$ws.Range("H136").Value = 4978.091
$ws.Range("I136").Value = 3969.875
$ws.Range("J136").Value = 7666.6665
$ws.Range("K136").Value = 11909.625
$ws.Range("L136").Value = 22999.9995
$ws.Range("M136").Value = -9359.625
$ws.Range("N136").Value = -28099.9995
$ws.Range("H140").Value = 27428.5
$ws.Range("J140").Value = 27428.5
$ws.Range("L140").Value = 27428.5
$ws.Range("N140").Value = -37788.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 9000
$ws.Range("I24").Value = 9000
$ws.Range("K24").Value = 9000
$ws.Range("M24").Value = -8770
$ws.Range("H29").Value = 45
$ws.Range("J29").Value = 45
$ws.Range("L29").Value = 45
$ws.Range("N29").Value = -625
$ws.Range("H62").Value = 98947.336
$ws.Range("I62").Value = 4341.2856
$ws.Range("K62").Value = 4341.2856
$ws.Range("M62").Value = -3717.2856
$ws.Range("H65").Value = 98947.336
$ws.Range("I65").Value = 4341.2856
$ws.Range("K65").Value = 21706.428
$ws.Range("M65").Value = -18586.428
$ws.Range("H81").Value = 3657.5625
$ws.Range("J81").Value = 9682.6
$ws.Range("L81").Value = 19365.2
$ws.Range("N81").Value = -21487.2
$ws.Range("H84").Value = 3657.5625
$ws.Range("J84").Value = 9682.6
$ws.Range("L84").Value = 96826
$ws.Range("N84").Value = -107434
$ws.Range("H122").Value = 2166.2144
$ws.Range("I122").Value = 2166.2144
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6498.6432
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4048.6432
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 355316.72
$ws.Range("I132").Value = 422887.12
$ws.Range("K132").Value = 1268661.36
$ws.Range("M132").Value = -1266131.36
$ws.Range("H136").Value = 3722.85
$ws.Range("I136").Value = 3391.8667
$ws.Range("K136").Value = 10175.6001
$ws.Range("M136").Value = -7625.6001
$ws.Range("H140").Value = 95999
$ws.Range("J140").Value = 95999
$ws.Range("L140").Value = 95999
$ws.Range("N140").Value = -106359
$ws.Range("H141").Value = 88198
$ws.Range("J141").Value = 88198
$ws.Range("L141").Value = 88198
$ws.Range("N141").Value = -98558

Write-Host "Applied Louisoix_Profits.xlsx updates across all sheets"